$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.712.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.29%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.893.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.66%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.85%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'311.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.76%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  -0.86%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4910"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.04%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -0.52%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -0.64%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.9119"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.38%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'20.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.04%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'1.933.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.34%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = "'TRON"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.07643"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.84%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.466"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.28%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'6.632"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.20%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'91.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.50%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -0.93%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000008777"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.14%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.06%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'27.902.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.59%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'14.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.18%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.119"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.13%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'2.138.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.40%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  -2.09%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  -2.27%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'1.860"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.43%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -1.11%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'2.154"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.95%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'115.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.59%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -2.11%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.08931"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.39%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  -3.94%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.229"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.60%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.7641"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.35%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'4.625"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.71%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = "'RenderToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'2.561"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -6.57%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'VeChain"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.02041"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.33%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -2.57%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -1.67%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.5462"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.39%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -0.81%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'6.894"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.53%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'8.554"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'113.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +7.47%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.1523"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.19%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'10.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.63%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.4792"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.21%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.9998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.91%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -2.27%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -1.57%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.06053"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.89%  "
$ws.Range("E51").Style = "Normal"
